$d = $word.ActiveDocument

# The document contains three inline pictures embedded in the headers/
# footers of the (single) section:
#   - Footer, type 1 (primary)  -> Pearson logo,  currently named "image1.png"
#   - Footer, type 2 (even pg)  -> Pearson logo,  currently named "image1.png"
#   - Header, type 2 (even pg)  -> BTEC logo,     currently named "image2.jpg"
#
# The edit simply renames the pictures (InlineShape.Name), swapping the
# numeric suffixes used by the two Pearson logo instances with the BTEC
# logo instance: image1.png -> image2.png (both Pearson copies) and
# image2.jpg -> image1.jpg (the BTEC copy). No other properties change.

foreach ($sec in $d.Sections) {

    # --- Footers -------------------------------------------------------
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.Name -eq "image1.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }

    # --- Headers ---------------------------------------------------------
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.Name -eq "image2.jpg") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }
}
